$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")

# Status: draft -> active
$ws1.Cells.Item(6, 2).Value = "active"

# Date: updated publish date
$ws1.Cells.Item(8, 2).Value = "2024-12-16T14:50:05-03:00"

# Case Sensitive: (blank/false) -> true
# Needs to land as literal text "true" (not a Boolean), matching how the
# rest of this metadata sheet stores true/false as plain strings.
# A leading apostrophe forces Excel to keep it as text, then we paste the
# formatting (not the value) from an untouched neighbour cell so the
# cell's style index is unaffected by the temporary quote-prefix flag.
$target = $ws1.Cells.Item(17, 2)
$formatDonor = $ws1.Cells.Item(19, 2)
$target.Value = "'true"
$formatDonor.Copy()
$target.PasteSpecial(-4122)
